# Update the workbook with the latest daily/monthly readings (atualizacao 16 nov 2020)

$wb = $excel.ActiveWorkbook

# --- Sheet "Mensal" (sheet1): append one new monthly summary row ---
$wsMensal = $wb.Worksheets.Item("Mensal")

$mensalRow = 14
$wsMensal.Cells.Item($mensalRow, 1).Value = 44150
$wsMensal.Cells.Item($mensalRow, 2).Value = 174.36
$wsMensal.Cells.Item($mensalRow, 3).Value = 208.88
$wsMensal.Cells.Item($mensalRow, 4).Value = -16.53

# copy the formatting (date number format / font / border) of the previous last row
$wsMensal.Range("A13").Copy()
$wsMensal.Range("A14").PasteSpecial(-4122) # xlPasteFormats
$wsMensal.Range("B13:D13").Copy()
$wsMensal.Range("B14:D14").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# --- Sheet "Diario" (sheet2): append fifteen new daily rows ---
$wsDiario = $wb.Worksheets.Item("Diario")

$dailyData = @(
    @(44136, 178.41, 208.88, -14.59),
    @(44137, 176.86, 208.88, -15.33),
    @(44138, 174.3,  208.88, -16.55),
    @(44139, 173.11, 208.88, -17.12),
    @(44140, 172.05, 208.88, -17.63),
    @(44141, 170.45, 208.88, -18.4),
    @(44142, 170.82, 208.88, -18.22),
    @(44143, 172.95, 208.88, -17.2),
    @(44144, 172.16, 208.88, -17.58),
    @(44145, 178.46, 208.88, -14.56),
    @(44146, 174.95, 208.88, -16.24),
    @(44147, 174.97, 208.88, -16.23),
    @(44148, 174.4,  208.88, -16.51),
    @(44149, 174.34, 208.88, -16.53),
    @(44150, 177.15, 208.88, -15.19)
)

$startRow = 368
for ($i = 0; $i -lt $dailyData.Count; $i++) {
    $r = $startRow + $i
    $vals = $dailyData[$i]
    $wsDiario.Cells.Item($r, 1).Value = $vals[0]
    $wsDiario.Cells.Item($r, 2).Value = $vals[1]
    $wsDiario.Cells.Item($r, 3).Value = $vals[2]
    $wsDiario.Cells.Item($r, 4).Value = $vals[3]
}

# copy formatting from the previous last row (367) down to the new rows
$wsDiario.Range("A367").Copy()
$wsDiario.Range("A368:A382").PasteSpecial(-4122) # xlPasteFormats
$wsDiario.Range("B367:D367").Copy()
$wsDiario.Range("B368:D382").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0
